# slit_plan_M75_BAO.xlsx - add Width-3/Width-4 and Weight-3/Weight-4 columns
#
# The sheet used to report only two slit widths/weights (Width-1/2,
# Weight-1/2) plus a combined "Extra Width Generated"/"Extra Weight"
# column. The plan now always produces up to four slit widths, so we add
# dedicated Width-3 / Width-4 columns (right after Width-2) and Weight-3 /
# Weight-4 columns (right after Weight-2), and recompute the two sample
# rows for the new upper bound (4 slits instead of 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the new columns -------------------------------
# Insert two blank columns right after "Width-2 (mm)" (column C) for the
# new Width-3 / Width-4 fields. This pushes everything from the old
# column D onward two slots to the right.
$ws.Columns("D:E").Insert()

# Insert two more blank columns right after "Weight-2 (kg)" (now column
# M) for the new Weight-3 / Weight-4 fields, ahead of "Extra Weight (kg)".
$ws.Columns("N:O").Insert()

# The original Width-3/Width-4/Weight-3/Weight-4 columns (which used to
# sit at the end of the table, after "Extra Weight") got shifted out to
# Q:T by the inserts above and are now duplicates of the new D:E/N:O
# columns - drop them.
$ws.Columns("Q:T").Delete()

# --- 2. Header row (comments / labels for the new columns) ----------
$ws.Range("D1").Value = "Width-3 (mm)"
$ws.Range("E1").Value = "Width-4 (mm)"
$ws.Range("N1").Value = "Weight-3 (kg)"
$ws.Range("O1").Value = "Weight-4 (kg)"

# --- 3. Recomputed data rows (upper bound raised from 2 to 4 slits) --
# Row 2 (plan id 0)
$ws.Range("B2").Value = 110
$ws.Range("C2").Value = 130
$ws.Range("D2").Value = 170
$ws.Range("E2").Value = 190
$ws.Range("F2").Value = 90
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 700
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 822.6491771451605
$ws.Range("L2").Value = 157.1428571428571
$ws.Range("M2").Value = 185.7142857142857
$ws.Range("N2").Value = 242.8571428571429
$ws.Range("O2").Value = 271.4285714285714
$ws.Range("P2").Value = 128.5714285714286

# Row 3 (plan id 1)
$ws.Range("B3").Value = 150
$ws.Range("C3").Value = 180
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 370
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 710
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 1216.593853524533
$ws.Range("L3").Value = 316.9014084507042
$ws.Range("M3").Value = 380.281690140845
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 781.6901408450702
